$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the date column (A2:A6) from 2025-11-29 to 2025-12-01
# Keep it stored as text (not an Excel date serial number). Temporarily
# mark the range as Text so the assignment isn't auto-parsed into a date
# serial, then restore the original (default) style so no formatting
# change is left behind.
$ws.Range("A2:A6").NumberFormat = "@"
$ws.Range("A2:A6").Value = "2025-12-01"
$ws.Range("A2:A6").Style = "Normal"

# Update the N column values (N2:N6) to the refreshed figure
$ws.Range("N2:N6").Value = 85.87246918135976
